# Commit: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# (Remove old account-statement (EC) worker rows and add new ones; also
# update the "old BD" summary figures.)
#
# The sheet previously listed 4 workers (rows 16-19) under the same NIT.
# Three of them (rows 17-19: HASMED MOISES CASTRO VITAL / DANIEL ANDRES
# CASTRO VITAL / MARNELLY ACOSTA POLANCO) are removed, and the single
# worker that remains (row 16) is replaced with the "new" worker
# (MARNELLY ACOSTA POLANCO / 1002280059), carrying the dues figures that
# used to belong to row 17 (56940 / 1423500). The summary fields (count
# of workers/periods, and total overdue value) are updated to match the
# now-single remaining worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 becomes the sole worker row, now holding the "new" identity and
# the dues values that used to sit in row 17.
$ws.Range("C16").Value = "1002280059"
$ws.Range("D16").Value = "MARNELLY ACOSTA POLANCO"
$ws.Range("E16").HorizontalAlignment = -4108  ; # xlCenter, matches the rest of the row's styling
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Drop the other three worker rows entirely.
$ws.Rows("17:19").Delete()

# Update the header summary: total overdue value and worker count.
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
